# "SAM TODO" workbook update — updated to do list
#
# 1. Insert a new "Not done" row right after the existing
#    "Implementation of IEC 61853 algorithms in C++" row (old row 64),
#    pushing the "Clean up the multiple load metrics..." row (old row 65)
#    and everything below it down by one.
# 2. Append two new "Future" rows just before the trailing "NGTD" row
#    (.zsam project file importer).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert new "Not done" row at row 65 ---------------------------
# Inserting at row 65 shifts the old row 65 ("Clean up the multiple load
# metrics reported to UI") and all following rows down by one, and the new
# blank row 65 inherits the formatting (fill/style) of the row above it.
$ws.Rows.Item(65).Insert()

# --- Step 2: insert two new "Future" rows before the trailing NGTD row -----
# After step 1 the NGTD row ("­.zsam project file importer") lives at row 81.
# Inserting twice at row 81 shifts it down to row 83 and leaves two blank
# rows (81, 82) that inherit the "Future" row formatting from above.
$ws.Rows.Item(81).Insert()
$ws.Rows.Item(81).Insert()

# Fill in the new cell text in the same order the new strings first appear
# in the shared-string table (resource-page summary, then outputs-browser
# group names, then the detailed SSC group-name fix) so the regenerated
# sharedStrings.xml keeps the same ordering as the target workbook.
$ws.Range("A82").Value = "Future"
$ws.Range("B82").Value = "Show annual summary statistics (irradiance, wind speed, ambient temp) on the resource page"

$ws.Range("A81").Value = "Future"
$ws.Range("B81").Value = "Fix output variable group names in SSC so that they show up pretty in outputs browser"

$ws.Range("A65").Value = "Not done"
$ws.Range("B65").Value = "Fix output variable group names in SSC: Flat Plate, PVWatts, Wind, CSP Trough & Towers, Res, Comm"

# --- Step 3: update the view state shown in the diff -----------------------
$ws.Application.ActiveWindow.ScrollRow = 59
$ws.Range("A66").Select()
